# Add a new data row (row 19) to Sheet1, mirroring the existing rows'
# formatting, with the same values/types the source workbook ends up with.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteValues  = -4163
$xlPasteFormats = -4122

# --- Plain text values: these don't look like numbers/booleans, so a
# normal .Value assignment keeps them as text without any type coercion. ---
$ws.Range("A19").Value = "TEST"
$ws.Range("B19").Value = "TEST"
$ws.Range("C19").Value = "6BF02F00"

# --- Numeric value ---
$ws.Range("D19").Value = 5

# --- Text values that LOOK like a number ("102") or a boolean ("TRUE").
# Assigning these directly via .Value would make Excel auto-convert them
# to a real number/boolean. Instead, build each one as a text formula
# (="102" / ="TRUE") in a scratch cell, convert that formula to a plain
# value in place (Copy + PasteSpecial values), then copy that already
# text-typed value into the destination cell. Finally, re-apply the
# destination's original number format/alignment (PasteSpecial formats)
# since the scratch cell's style differs from the target row's style. ---

# E19 <- "102"
$ws.Range("H1").Formula = "=""102"""
$ws.Range("H1").Copy()
$ws.Range("H1").PasteSpecial($xlPasteValues)
$ws.Range("H1").Copy()
$ws.Range("E19").PasteSpecial($xlPasteValues)
$ws.Range("H1").Clear()
$ws.Range("E18").Copy()
$ws.Range("E19").PasteSpecial($xlPasteFormats)

# F19 <- "TRUE"
$ws.Range("H1").Formula = "=""TRUE"""
$ws.Range("H1").Copy()
$ws.Range("H1").PasteSpecial($xlPasteValues)
$ws.Range("H1").Copy()
$ws.Range("F19").PasteSpecial($xlPasteValues)
$ws.Range("H1").Clear()
$ws.Range("F18").Copy()
$ws.Range("F19").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0
